# Updated test_2ci and test_NAs to run properly
$wb = $excel.ActiveWorkbook

# --- testdata_Mean (sheet1): remove the final "Area2 / NA" row and the
#     trailing blank row so test_2ci only sees complete data ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A29:B29").ClearContents() | Out-Null
$ws1.Rows.Item(30).Delete() | Out-Null
$ws1.Rows.Item(29).Select() | Out-Null

# --- testdata_Mean_results_NA (new sheet3): duplicate of
#     testdata_Mean_results but with Area1's values replaced by "NA" so
#     test_NAs has a fixture to run against ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "testdata_Mean_results_NA"

$ws3.Range("B2").Value = "NA"
$ws3.Range("D2").Value = "NA"
$ws3.Range("E2").Value = "NA"
$ws3.Range("F2").Value = "NA"
$ws3.Range("G2").Value = "NA"
$ws3.Range("H2").Value = "NA"
$ws3.Range("I2").Value = "NA"

$ws3.Range("A2").Select() | Out-Null
$ws3.Activate() | Out-Null
